# "more tests and fixes"
#
# Changes applied to tests/test_engines/engines.xlsx:
#  1. Sheet1!A1 text changes from "Bye xlwings!" to "a" (shared string update).
#  2. Sheet1's selection/active cell moves from B15 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure Sheet1 is the active sheet (it already is tabSelected in the
# original file) so the selection below is recorded against it.
$ws.Activate()

# 1) Update the cell text.
$ws.Range("A1").Value = "a"

# 2) Move the selection to A2 (was B15).
$ws.Range("A2").Select() | Out-Null
